$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# for every data row (rows 2 through 308).
$ws.Range("C2:C308").Value = 45189
